$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 266, shifting rows 266:372 down to 267:373
$ws.Rows.Item(266).Insert()

# Populate the new row 266 with the new data record
$ws.Cells.Item(266, 1).Value = 11
$ws.Cells.Item(266, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(266, 3).Value = "Bíobío"
$ws.Cells.Item(266, 4).Value = 45027
$ws.Cells.Item(266, 5).Value = 8
$ws.Cells.Item(266, 6).Value = 100112045
$ws.Cells.Item(266, 7).Value = "Zapallo"
$ws.Cells.Item(266, 8).Value = "Camote"
$ws.Cells.Item(266, 9).Value = "1a (cosecha)"
$ws.Cells.Item(266, 10).Value = 900
$ws.Cells.Item(266, 11).Value = 300
$ws.Cells.Item(266, 12).Value = 350
$ws.Cells.Item(266, 13).Value = 322
$ws.Cells.Item(266, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(266, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(266, 16).Value = 322
$ws.Cells.Item(266, 17).Value = 1
$ws.Cells.Item(266, 18).Value = "Hortaliza"

# Apply the same style as column D (date) cells to D266
$ws.Cells.Item(266, 4).NumberFormat = $ws.Cells.Item(267, 4).NumberFormat
